# "Wall day 1 part 2" -- append two Tag-log rows (actions 101/102) to the
# game log table on Sheet1, and grow Table1 to accommodate future rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- expand the table first, like Excel does when you type past its edge ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G110"))

# --- row 99: Riya Bhangu tags Mya Wiggins ---
$ws.Cells.Item(99, 1).Value = 101
$ws.Cells.Item(99, 2).Value = "Tag"
$ws.Cells.Item(99, 3).Value = "Riya Bhangu"
$ws.Cells.Item(99, 4).Value = "Mya Wiggins"
$ws.Cells.Item(99, 5).Value = 45755
$ws.Cells.Item(99, 6).Value = 0.48541666666666666

# --- row 100: Haley Jones tags Ellie Milligan ---
$ws.Cells.Item(100, 1).Value = 102
$ws.Cells.Item(100, 2).Value = "Tag"
$ws.Cells.Item(100, 3).Value = "Haley Jones"
$ws.Cells.Item(100, 4).Value = "Ellie Milligan"
$ws.Cells.Item(100, 5).Value = 45756
$ws.Cells.Item(100, 6).Value = 0.46250000000000002

# --- column G: continue the Unix-time helper formula, matching the
#     look/format of the existing rows (row 98's style) ---
$fmtSrc = $ws.Cells.Item(98, 7)
for ($r = 99; $r -le 100; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Formula = "=((E$r+F$r)-DATE(1970,1,1))*86400"
    $cell.NumberFormat = $fmtSrc.NumberFormat
    $cell.Font.Name = $fmtSrc.Font.Name
    $cell.Font.Size = $fmtSrc.Font.Size
    $cell.Font.Color = $fmtSrc.Font.Color
    $cell.Borders.Item(9).LineStyle = $fmtSrc.Borders.Item(9).LineStyle
    $cell.HorizontalAlignment = $fmtSrc.HorizontalAlignment
}

# --- match the scrolled view / active selection left by the edit ---
$ws.Range("A101").Select() | Out-Null
